$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Formula = "=ROW()*COLUMN()"
$ws.Range("A2").Formula = "=A1*10"
$ws.Range("B3").Formula = "=A2+B1"

$ws.Range("B3").Select() | Out-Null
